$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last 7 rows of the "Testing Table" (rows 37-43) never had their
# "Test Number" column filled in, and picked up a slightly different
# (but visually identical) cell style along the way. Finish numbering the
# sequence, continuing on from row 36 (test #34).
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
$ws.Range("A40").Value = 38
$ws.Range("A41").Value = 39
$ws.Range("A42").Value = 40
$ws.Range("A43").Value = 41

# Normalize the formatting of those rows to match the rest of the table by
# copying the formats down from row 36 (the last "clean" row).
$ws.Range("A36:H36").Copy()
[void]$ws.Range("A37:H43").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Scroll back to the left and select the whole table body, matching where
# the author left the cursor after finishing this pass (still need to
# update the wind values next).
[void]$ws.Range("A2:H43").Select()
